$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '57.110.40'
$ws.Range("E2").Value = '  +7.16%  '

# Row 3
$ws.Range("D3").Value = '3.240.67'
$ws.Range("E3").Value = '  +2.82%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").Value = '''395.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.64%  '

# Row 6
$ws.Range("D6").Value = '''107.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.39%  '

# Row 7
$ws.Range("D7").Value = '3.239.76'
$ws.Range("E7").Value = '  +2.96%  '

# Row 8
$ws.Range("D8").Value = '''0.569'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.58%  '

# Row 9
$ws.Range("E9").Value = '  +0.02%  '

# Row 10
$ws.Range("E10").Value = '  +1.23%  '

# Row 11
$ws.Range("D11").Value = '''38.94'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.29%  '

# Row 12
$ws.Range("D12").Value = '''0.0971'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +11.58%  '

# Row 13
$ws.Range("E13").Value = '  +1.70%  '

# Row 14
$ws.Range("D14").Value = '3.753.64'
$ws.Range("E14").Value = '  +2.98%  '

# Row 15
$ws.Range("D15").Value = '''8.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.54%  '

# Row 16
$ws.Range("D16").Value = '''18.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.47%  '

# Row 17
$ws.Range("D17").Value = '3.245.56'
$ws.Range("E17").Value = '  +3.06%  '

# Row 18
$ws.Range("E18").Value = '  -2.27%  '

# Row 19
$ws.Range("D19").Value = '''11.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.54%  '

# Row 20
$ws.Range("D20").Value = '56.954.01'
$ws.Range("E20").Value = '  +6.90%  '

# Row 21
$ws.Range("E21").Value = '  +1.02%  '

# Row 22
$ws.Range("E22").Value = '  +8.93%  '

# Row 23
$ws.Range("D23").Value = '''13.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.43%  '

# Row 24
$ws.Range("D24").Value = '''296.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.56%  '

# Row 25
$ws.Range("D25").Value = '''73.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.41%  '

# Row 26
$ws.Range("D26").Value = '''3.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.14%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''27.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.57%  '

# Row 28
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").Value = '''7.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.33%  '

# Row 29
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '''7.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.09%  '

# Row 30
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '''0.169'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.01%  '

# Row 31
$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").Value = '''1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '

# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.108'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.64%  '

# Row 33
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '''11.03'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.41%  '

# Row 34
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = '''37.30'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.19%  '

# Row 35
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = '''0.0484'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.91%  '

# Row 36
$ws.Range("B36").Value = 'Toncoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D36").Value = '''2.11'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.75%  '

# Row 37
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '''51.69'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.80%  '

# Row 38
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '''3.53'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.41%  '

# Row 39
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '''0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.14%  '

# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''3.03'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.99%  '

# Row 41
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").Value = '''134.48'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.30%  '

# Row 42
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '''0.120'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.33%  '

# Row 43
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '''1.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.44%  '

# Row 44
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").Value = '''3.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.14%  '

# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = '''16.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.94%  '

# Row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = '''0.281'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.55%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''21.85'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.38%  '

# Row 48
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.148.09'
$ws.Range("E48").Value = '  +2.72%  '

# Row 49
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = '''2.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.78%  '

# Row 50
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '''2.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +23.55%  '

# Row 51
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = '''2.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.43%  '
